# Remove the trailing "subscribe to the RSS feed" paragraph, the
# horizontal-rule paragraph that follows it, and the final blank
# paragraph at the very end of the document body.
$d = $word.ActiveDocument

# Locate the start of the paragraph that begins the text to be removed.
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "This post is a first look at a dataset of job openings",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $startPara = $searchRange.Paragraphs.Item(1)
    $deleteStart = $startPara.Range.Start
    $deleteEnd = $d.Content.End

    $target = $d.Range($deleteStart, $deleteEnd)
    $target.Delete()
}
